$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cluster labels for rows 2-5 (columns A-D): Sending, Ligand, Receptor, Target
$data = @(
    @("FAPs","Plg","F3","FAPs", 1, 0.3333333333333333, 0.465759, 1.397277, 0.9660838355812051, 0.9660838355812051, 3, 1, 173.637756, 520.913268, 0.9875738673498291, 0.9875738673498291, 80.87334759680401, 727.8601283712361, 0.9540791496890872, 0.9540791496890872),
    @("FAPs","Plg","F3","MuSCs", 1, 0.3333333333333333, 0.465759, 1.397277, 0.9660838355812051, 0.9660838355812051, 3, 1, 2.184794333333333, 6.554383, 0.0124261326501708, 0.0124261326501708, 1.017587623899, 9.158288615091001, 0.01200468589211786, 0.01200468589211786),
    @("MuSCs","Plg","F3","FAPs", 1, 0.3333333333333333, 0.01635133333333333, 0.049054, 0.03391616441879487, 0.03391616441879487, 3, 1, 173.637756, 520.913268, 0.9875738673498291, 0.9875738673498291, 2.839208827608, 25.552879448472, 0.03349471766074191, 0.03349471766074191),
    @("MuSCs","Plg","F3","MuSCs", 1, 0.3333333333333333, 0.01635133333333333, 0.049054, 0.03391616441879487, 0.03391616441879487, 3, 1, 2.184794333333333, 6.554383, 0.0124261326501708, 0.0124261326501708, 0.03572430040911111, 0.321518703682, 0.0004214467580529482, 0.0004214467580529482)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}

$wb.Save()
